$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 - Suecia: only "Muertes hoy" (F) changed
$ws.Range("F24").Value = 523

# Row 26 - Corea del Sur
$ws.Range("B26").Value = 10653
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 7937
$ws.Range("E26").Value = 2484
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 232

# Row 37 - Australia
$ws.Range("D37").Value = 4132
$ws.Range("E37").Value = 2361

# Row 48/49 - Panama moves above Republica Dominicana (Panama got new, higher totals)
$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 4210
$ws.Range("C48").Value = 194
$ws.Range("D48").Value = 122
$ws.Range("E48").Value = 3972
$ws.Range("F48").Value = 96
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 116

$ws.Range("A49").Value = "Republica Dominicana"
$ws.Range("B49").Value = 4126
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 268
$ws.Range("E49").Value = 3658
$ws.Range("F49").Value = 146
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 200

# Row 70 - Nueva Zelanda
$ws.Range("B70").Value = 1422
$ws.Range("C70").Value = 13
$ws.Range("D70").Value = 867
$ws.Range("E70").Value = 544
$ws.Range("F70").Value = 3

# Row 120/121 - Venezuela moves above Guatemala (Venezuela got new, higher totals)
$ws.Range("A120").Value = "Venezuela"
$ws.Range("B120").Value = 227
$ws.Range("C120").Value = 23
$ws.Range("D120").Value = 111
$ws.Range("E120").Value = 107
$ws.Range("F120").Value = 4
$ws.Range("H120").Value = 9

$ws.Range("A121").Value = "Guatemala"
$ws.Range("B121").Value = 214
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 21
$ws.Range("E121").Value = 186
$ws.Range("F121").Value = 3
$ws.Range("H121").Value = 7

# Row 122 - Paraguay
$ws.Range("B122").Value = 202
$ws.Range("C122").Value = 3
$ws.Range("D122").Value = 35
$ws.Range("E122").Value = 159
